$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.037.44'

# Row 3
$ws.Range("D3").Value = '1.833.18'
$ws.Range("E3").Value = '  -0.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.29'
$ws.Range("E5").Value = '  -0.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6227'
$ws.Range("E6").Value = '  -6.26%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.21%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07566'
$ws.Range("E8").Value = '  +1.57%  '

# Row 9
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.59'
$ws.Range("E9").Value = '  -0.66%  '

# Row 10
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2914'
$ws.Range("E10").Value = '  -1.60%  '

# Row 11
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.59'
$ws.Range("E11").Value = '  -3.25%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07741'
$ws.Range("E12").Value = '  -0.17%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.840.35'
$ws.Range("E13").Value = '  -0.39%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.945'
$ws.Range("E14").Value = '  -1.57%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6645'
$ws.Range("E15").Value = '  -1.37%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009976'
$ws.Range("E16").Value = '  +14.18%  '

# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '82.65'
$ws.Range("E17").Value = '  -0.59%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.037'
$ws.Range("E18").Value = '  -2.14%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '29.054.44'
$ws.Range("E19").Value = '  -0.53%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '225.66'
$ws.Range("E20").Value = '  -0.56%  '

# Row 21
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.34'
$ws.Range("E21").Value = '  -1.73%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.12%  '

# Row 23
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.170'
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.28%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.39'
$ws.Range("E25").Value = '  +0.22%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1375'
$ws.Range("E26").Value = '  -1.96%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.440'
$ws.Range("E27").Value = '  -2.21%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.89'
$ws.Range("E28").Value = '  -0.94%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.488'
$ws.Range("E29").Value = '  -1.67%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.069'
$ws.Range("E30").Value = '  -1.69%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.029'
$ws.Range("E31").Value = '  -0.58%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.201'
$ws.Range("E32").Value = '  +0.61%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05201'
$ws.Range("E33").Value = '  -2.56%  '

# Row 34
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.847'
$ws.Range("E34").Value = '  -0.73%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7382'
$ws.Range("E35").Value = '  -1.15%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.139'
$ws.Range("E36").Value = '  -1.78%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.700'
$ws.Range("E37").Value = '  +1.73%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.246.92'
$ws.Range("E38").Value = '  -4.88%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.764'
$ws.Range("E39").Value = '  +0.33%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01783'
$ws.Range("E40").Value = '  -0.77%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.371'
$ws.Range("E41").Value = '  -0.68%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8956'
$ws.Range("E42").Value = '  -1.17%  '

# Row 43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.20%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.50'
$ws.Range("E44").Value = '  -1.91%  '

# Row 45
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.983.05'
$ws.Range("E45").Value = '  -0.33%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000124'
$ws.Range("E46").Value = '  +0.01%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.06'
$ws.Range("E47").Value = '  -1.56%  '

# Row 48
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5115'
$ws.Range("E48").Value = '  -0.33%  '

# Row 49
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4017'
$ws.Range("E49").Value = '  -0.12%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.843'
$ws.Range("E50").Value = '  +0.32%  '

# Row 51
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.643'
$ws.Range("E51").Value = '  -6.22%  '
